# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-03-16 Sunday", $true, $false, $false, $false, `
    $false, $true, 1, $false, "2025-03-17 Monday", 2)

# Update the division problems in the table, row by row / cell by cell,
# to avoid ambiguity with duplicate values (e.g. "61÷4=" appears twice).
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("88÷6=", "36÷8=", "77÷4=", "18÷2=", "84÷5=")
    5  = @("90÷2=", "24÷7=", "41÷8=", "81÷4=", "92÷2=")
    9  = @("24÷3=", "43÷9=", "86÷7=", "43÷8=", "95÷3=")
    13 = @("76÷3=", "19÷6=", "45÷4=", "48÷3=", "28÷6=")
    17 = @("59÷3=", "85÷9=", "73÷2=", "11÷8=", "51÷5=")
}

foreach ($row in $newValues.Keys) {
    $values = $newValues[$row]
    for ($c = 1; $c -le $values.Length; $c++) {
        $t.Cell($row, $c).Range.Text = $values[$c - 1]
    }
}
